$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header D1 from "IsAvailable" to "Availability"
$ws.Range("D1").Value = "Availability"

# Change D2:D5 from boolean TRUE to numeric 0
$ws.Range("D2:D5").Value = 0

# Update sheet view: clear topLeftCell, change selection to I9
$ws.Range("I9").Select()
